$wb = $excel.ActiveWorkbook
$s = $wb.Styles.Item(1)
$name = "Norm" + [char]0x00E1 + "ln" + [char]0x00ED
$s.Name = $name
